$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.454.50"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "3.982.79"
$ws.Range("E3").Value = "  +4.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.37%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.737"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.171"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000368"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "4.604.82"
$ws.Range("E13").Value = "  +4.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "4.003.10"
$ws.Range("E16").Value = "  +4.96%  "
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.67%  "
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "68.492.19"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.70%  "
$ws.Range("E22").Value = "  +4.72%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("E27").Value = "  +5.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "733.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("E32").Value = "  +3.96%  "
$ws.Range("D33").Value = "0.0₃0906"
$ws.Range("E33").Value = "  +34.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.83%  "
$ws.Range("E36").Value = "  -3.42%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0476"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  +5.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.31%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.142"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.06%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.12%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.66%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.338"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.50%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.89%  "
$ws.Range("E49").Value = "  +2.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.08%  "
